$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix capitalization of class names in column A (rows 2, 4, 8)
$ws.Range("A2").Value = "mdaTextHomePage"
$ws.Range("A8").Value = "pageTitleNewTab"
$ws.Range("A4").Value = "mdaTitle"

# Update selected / active cell on the sheet
$ws.Range("A4").Select()
